$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/10/2026"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 9118.940000000001
$ws.Cells.Item($row, 3).Value = 0.2470222779281854
$ws.Cells.Item($row, 4).Value = 0.7529777220718146
$ws.Cells.Item($row, 5).Value = -333.72
$ws.Cells.Item($row, 6).Value = -39.97
$ws.Cells.Item($row, 7).Value = -24025.4
$ws.Cells.Item($row, 8).Value = -77.77
$ws.Cells.Item($row, 9).Value = -1100.27
$ws.Cells.Item($row, 10).Value = -32.82
$ws.Cells.Item($row, 11).Value = -25125.67
$ws.Cells.Item($row, 12).Value = -73.37
